$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewSynonym")
$ws.Columns("D").Insert()
$ws.Cells.Item(1, 4).Value = "predicate"
